$wb = $excel.ActiveWorkbook

# --- ip_address_list: insert a new row at the top, shift rows down,
#     move 'Domac' row, append extra text, flip a flag ---
$ws1 = $wb.Worksheets.Item("ip_address_list")
$ws1.Rows.Item(1).Insert()

# row 1: 'xz'
$ws1.Cells.Item(1,1).Value = "xz"
$ws1.Cells.Item(1,2).Value = "192.168.000.000"
$ws1.Cells.Item(1,3).Value = "255.255.255.0"
$ws1.Cells.Item(1,4).Value = ""
$ws1.Cells.Item(1,5).Value = "'1"

# row 2: '474 B_Austin'
$ws1.Cells.Item(2,1).Value = "474 B_Austin"
$ws1.Cells.Item(2,2).Value = "10.96.205.175"
$ws1.Cells.Item(2,3).Value = "255.255.255.0"
$ws1.Cells.Item(2,4).Value = "PC:`t10.96.aoj"
$ws1.Cells.Item(2,5).Value = "'1"

# row 3: '529_Witte55'
$ws1.Cells.Item(3,1).Value = "529_Witte55"
$ws1.Cells.Item(3,2).Value = "192.168.0.240"
$ws1.Cells.Item(3,3).Value = "255.255.255.0"
$ws1.Cells.Item(3,4).Value = "P"
$ws1.Cells.Item(3,5).Value = "'0"

# row 4: '474 B_Austin (1)'
$ws1.Cells.Item(4,1).Value = "474 B_Austin (1)"
$ws1.Cells.Item(4,2).Value = "10.96.205.175"
$ws1.Cells.Item(4,3).Value = "255.255.255.0"
$ws1.Cells.Item(4,4).Value = "PC:`t10.96.205.175`nNAS:`t10.96.205.166`nFH:`t10.96.205.154`n`t10.96.20`n-----------------------------------------`nuser:JHV_Vision, omron `nPass:*Jhv2708`n---------------------------------------`nFortiClient Austin: `nPass:`n1Pm#J@PFIkzM&Q@i `nUVt1@Ex2p78kxp30atD7we@!qGK`ndfa"
$ws1.Cells.Item(4,5).Value = "'1"

# row 5: '474 B_Austin (2)'
$ws1.Cells.Item(5,1).Value = "474 B_Austin (2)"
$ws1.Cells.Item(5,2).Value = "10.96.205.175"
$ws1.Cells.Item(5,3).Value = "255.255.255.0"
$ws1.Cells.Item(5,4).Value = "10.96.205.1`nNAS:`t10.96.205.166`nFH:`t10.96.205.154`n`t10.96.20`n-----------------------------------------`nuser:JHV_Vision, omron llllllllllllll`nPass:*Jhv2708`n---------------------------------------`nFortiClient Austin: `nPass:`n1Pm#J@PFIkzM&Q@i `nUVt1@Ex2p78kxp30atD7we@!qGK"
$ws1.Cells.Item(5,5).Value = "'0"

# row 6: '474 B_Austin (2) (1)'
$ws1.Cells.Item(6,1).Value = "474 B_Austin (2) (1)"
$ws1.Cells.Item(6,2).Value = "10.96.205.175"
$ws1.Cells.Item(6,3).Value = "255.255.255.0"
$ws1.Cells.Item(6,4).Value = "10.96.205.1`nNAS:`t10.96.205.166`nFH:`t10.96.205.154`n`t10.96.20`n-----------------------------------------`nuser:JHV_Vision, omron llllllllllllll`nPass:*Jhv2708`n---------------------------------------`nFortiClient Austin: `nPass:`n1Pm#J@PFIkzM&Q@i `nUVt1@Ex2p78kxp30atD7we@!qGK`naf"
$ws1.Cells.Item(6,5).Value = "'1"

# row 7: '529_Witte'
$ws1.Cells.Item(7,1).Value = "529_Witte"
$ws1.Cells.Item(7,2).Value = "192.168.0.240"
$ws1.Cells.Item(7,3).Value = "255.255.255.0"
$ws1.Cells.Item(7,4).Value = "PC"
$ws1.Cells.Item(7,5).Value = "'0"

# row 8: 'Domac'
$ws1.Cells.Item(8,1).Value = "Domac"
$ws1.Cells.Item(8,2).Value = "192.168.1.13"
$ws1.Cells.Item(8,3).Value = "255.255.255.0"
$ws1.Cells.Item(8,4).Value = ""
$ws1.Cells.Item(8,5).Value = "'0"

# row 9: '47'
$ws1.Cells.Item(9,1).Value = "'47"
$ws1.Cells.Item(9,2).Value = "10.96.205.175"
$ws1.Cells.Item(9,3).Value = "255.255.255.0"
$ws1.Cells.Item(9,4).Value = "PC:`t10.96.205.`nasdf"
$ws1.Cells.Item(9,5).Value = "'1"

# row 10: '518_Val'
$ws1.Cells.Item(10,1).Value = "518_Val"
$ws1.Cells.Item(10,2).Value = "192.168.208.242"
$ws1.Cells.Item(10,3).Value = "255.255.255.0"
$ws1.Cells.Item(10,4).Value = ""
$ws1.Cells.Item(10,5).Value = "'0"

# row 11: '518_Valeo II'
$ws1.Cells.Item(11,1).Value = "518_Valeo II"
$ws1.Cells.Item(11,2).Value = "192.168.1.243"
$ws1.Cells.Item(11,3).Value = "255.255.255.0"
$ws1.Cells.Item(11,4).Value = "ssssssss"
$ws1.Cells.Item(11,5).Value = "'0"

# row 12: '527_'
$ws1.Cells.Item(12,1).Value = "527_"
$ws1.Cells.Item(12,2).Value = "10.101.28.176"
$ws1.Cells.Item(12,3).Value = "255.255.255.0"
$ws1.Cells.Item(12,4).Value = "PC:`t10.96.20"
$ws1.Cells.Item(12,5).Value = "'0"

# row 13: 'Dom'
$ws1.Cells.Item(13,1).Value = "Dom"
$ws1.Cells.Item(13,2).Value = "192.168.1.131"
$ws1.Cells.Item(13,3).Value = "255.255.255.0"
$ws1.Cells.Item(13,4).Value = ""
$ws1.Cells.Item(13,5).Value = "'1"

# row 14: '474 B_A'
$ws1.Cells.Item(14,1).Value = "474 B_A"
$ws1.Cells.Item(14,2).Value = "10.96.205.175"
$ws1.Cells.Item(14,3).Value = "255.255.255.0"
$ws1.Cells.Item(14,4).Value = "dfddddddddddddddddd`nadf`nafd`nafsdfaadfs"
$ws1.Cells.Item(14,5).Value = "'0"

# --- disk_list: move 515_ZF row to the end, append extra text ---
$ws3 = $wb.Worksheets.Item("disk_list")

# row 1: '514_Teleflex'
$ws3.Cells.Item(1,1).Value = "514_Teleflex"
$ws3.Cells.Item(1,2).Value = "T"
$ws3.Cells.Item(1,3).Value = "\\192.168.14.245\Data\Kamery"
$ws3.Cells.Item(1,4).Value = "Vision"
$ws3.Cells.Item(1,5).Value = "*Jhv2708"
$ws3.Cells.Item(1,6).Value = ""

# row 2: 'Domaci Nas'
$ws3.Cells.Item(2,1).Value = "Domaci Nas"
$ws3.Cells.Item(2,2).Value = "S"
$ws3.Cells.Item(2,3).Value = "\\192.168.1.20\Data"
$ws3.Cells.Item(2,4).Value = ""
$ws3.Cells.Item(2,5).Value = ""
$ws3.Cells.Item(2,6).Value = ""

# row 3: '518_Valeo II'
$ws3.Cells.Item(3,1).Value = "518_Valeo II"
$ws3.Cells.Item(3,2).Value = "V"
$ws3.Cells.Item(3,3).Value = "\\192.168.1.10\10_vision"
$ws3.Cells.Item(3,4).Value = "jhv_vision"
$ws3.Cells.Item(3,5).Value = "Jhv*2708"
$ws3.Cells.Item(3,6).Value = "Druha sít, ixon"

# row 4: '518_Valeo'
$ws3.Cells.Item(4,1).Value = "518_Valeo"
$ws3.Cells.Item(4,2).Value = "V"
$ws3.Cells.Item(4,3).Value = "\\192.168.208.200\10_vision"
$ws3.Cells.Item(4,4).Value = "jhv_vision"
$ws3.Cells.Item(4,5).Value = "Jhv*2708"
$ws3.Cells.Item(4,6).Value = "první sít, ixon`n\\192.168.208.200\10_vision`nsadf"

# row 5: '474_B Austin'
$ws3.Cells.Item(5,1).Value = "474_B Austin"
$ws3.Cells.Item(5,2).Value = "P"
$ws3.Cells.Item(5,3).Value = "\\10.96.205.166\DATA"
$ws3.Cells.Item(5,4).Value = "jhv_vision"
$ws3.Cells.Item(5,5).Value = "*Jhv2708"
$ws3.Cells.Item(5,6).Value = "10.96.205.166`t`nVisionNas_474B`t`n`t`t`t`t`t`tuser:JHV_Vision, omron `nPass:*Jhv2708`nadf"

# row 6: '515_ZF'
$ws3.Cells.Item(6,1).Value = "515_ZF"
$ws3.Cells.Item(6,2).Value = "Z"
$ws3.Cells.Item(6,3).Value = "\\10.9.250.100\08_Project_ZF_515\kamery"
$ws3.Cells.Item(6,4).Value = "jhvadmin"
$ws3.Cells.Item(6,5).Value = "jhvadm1n"
$ws3.Cells.Item(6,6).Value = ""

# --- projects_bin2: move 511_Teleflex row up from row 3 to row 2 ---
$ws4 = $wb.Worksheets.Item("projects_bin2")
$ws4.Cells.Item(2,1).Value = "511_Teleflex"
$ws4.Cells.Item(2,2).Value = "192.168.1.242"
$ws4.Cells.Item(2,3).Value = "255.255.255.0"
$ws4.Cells.Item(2,4).Value = "Teleflex d"
$ws4.Cells.Item(2,5).Value = $false
$ws4.Cells.Item(3,1).Value = ""
$ws4.Cells.Item(3,2).Value = ""
$ws4.Cells.Item(3,3).Value = ""
$ws4.Cells.Item(3,4).Value = ""
$ws4.Cells.Item(3,5).Value = ""

